$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '30.539.99'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9993'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '236.77'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -3.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9989'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4873'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -1.85%  '
$ws.Range('E8').Value = '  -2.16%  '
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.871.25'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -1.15%  '
$ws.Range('E11').Value = '  -2.74%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07221'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -1.60%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '89.46'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.001'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.95%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.6537'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -3.08%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '30.484.40'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000007835'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.9989'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('E19').Value = '  -2.01%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '2.113.01'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.02%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.9954'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.48%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '213.78'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +19.34%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.726'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -2.88%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.128'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.10%  '
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '156.05'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +1.27%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.02'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.42%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.827'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -5.22%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.410'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +1.50%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.262'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -1.62%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.09037'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +1.19%  '
$ws.Range('E32').Value = '  -2.74%  '
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7235'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('E35').Value = '  -4.87%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.686'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.49%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.01812'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.31%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.657'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.61%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.9193'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  -5.73%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.4403'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +0.95%  '
$ws.Range('E42').Value = '  -1.46%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.734'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.9941'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.67%  '
$ws.Range('E45').Value = '  -2.01%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.339'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -4.05%  '
$ws.Range('E47').Value = '  +3.09%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.05823'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.38%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '8.599'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +1.16%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.408'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '33.16'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.63%  '
